$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 14:49"

# --- Updated statistics for existing countries (values only, no row moves) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2683301
$ws.Range("C4").Value = 1490
$ws.Range("D4").Value = 1122591
$ws.Range("E4").Value = 1431891
$ws.Range("G4").Value = 36
$ws.Range("H4").Value = 128819

# Row 7 - India
$ws.Range("B7").Value = 568536
$ws.Range("C7").Value = 1000
$ws.Range("D7").Value = 335915
$ws.Range("E7").Value = 215702

# Row 29 - Bielorrusia
$ws.Range("B29").Value = 62118
$ws.Range("C29").Value = 328
$ws.Range("D29").Value = 46054
$ws.Range("E29").Value = 15672
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 392

# Row 33 - Paises Bajos
$ws.Range("B33").Value = 50273
$ws.Range("C33").Value = 50
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 6113

# Row 34 - Emiratos Arabes Unidos
$ws.Range("B34").Value = 48667
$ws.Range("C34").Value = 421
$ws.Range("D34").Value = 37566
$ws.Range("E34").Value = 10786
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 315

# Row 36 - Kuwait
$ws.Range("B36").Value = 46195
$ws.Range("C36").Value = 671
$ws.Range("D36").Value = 37030
$ws.Range("E36").Value = 8811
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 354

# Row 47 - Afganistan
$ws.Range("B47").Value = 31517
$ws.Range("C47").Value = 279
$ws.Range("D47").Value = 14131
$ws.Range("E47").Value = 16640
$ws.Range("G47").Value = 13
$ws.Range("H47").Value = 746

# Row 57 - Austria
$ws.Range("B57").Value = 17766
$ws.Range("C57").Value = 43
$ws.Range("D57").Value = 16478
$ws.Range("E57").Value = 583
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 705

# Row 66 - Dinamarca
$ws.Range("B66").Value = 12768
$ws.Range("C66").Value = 17
$ws.Range("D66").Value = 11649
$ws.Range("E66").Value = 514

# Row 74 - Uzbekistan
$ws.Range("B74").Value = 8385
$ws.Range("C74").Value = 163
$ws.Range("E74").Value = 2791

# Row 77 - Consejo Danes para los Refugiados
$ws.Range("B77").Value = 7039
$ws.Range("C77").Value = 100
$ws.Range("D77").Value = 1426
$ws.Range("E77").Value = 5443
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 170

# Row 91 - Bosnia y Herzegovina
$ws.Range("B91").Value = 4453
$ws.Range("C91").Value = 128
$ws.Range("D91").Value = 2402
$ws.Range("E91").Value = 1865
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 186

# Row 101 - Croacia
$ws.Range("B101").Value = 2777
$ws.Range("C101").Value = 52
$ws.Range("E101").Value = 515

# Row 106 - Maldivas
$ws.Range("E106").Value = 401
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 9

# Row 115 - Islandia
$ws.Range("B115").Value = 1842
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 1819
$ws.Range("E115").Value = 13

# Row 125 - Hong Kong
$ws.Range("B125").Value = 1206
$ws.Range("C125").Value = 2
$ws.Range("D125").Value = 1107

# --- Reordering of three country pairs/triples in the shared-string list ---
# The underlying country-name cells keep the same shared-string slot but the
# slot's text (and therefore the row that displays it) shifts, while the
# numeric statistics for that physical row stay attached to the row number.
# Net visible effect: row 162/163/164 now read Siria/Angola/Comoras (with
# data shifted accordingly), row 203/204 now read Laos/Santa Lucia, and row
# 209/210 now read Islas Malvinas/Groenlandia.

# Rows 162-164 (previously Angola, Comoras, Siria)
$ws.Range("A162").Value = "Siria"
$ws.Range("B162").Value = 279
$ws.Range("C162").Value = 10
$ws.Range("D162").Value = 105
$ws.Range("E162").Value = 165
$ws.Range("H162").Value = 9

$ws.Range("A163").Value = "Angola"
$ws.Range("B163").Value = 276
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 93
$ws.Range("E163").Value = 172
$ws.Range("H163").Value = 11

$ws.Range("A164").Value = "Comoras"
$ws.Range("B164").Value = 272
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 161
$ws.Range("E164").Value = 104
$ws.Range("H164").Value = 7

# Rows 203-204 (previously Santa Lucia, Laos) - values unaffected, only text swap
$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"

# Rows 209-210 (previously Groenlandia, Islas Malvinas) - values unaffected, only text swap
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
